$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new question row (row 39)
$ws.Range("A39").Value = 98
$ws.Range("C39").Value = "Validate Binary Search Tree"

# Update selection to match the saved view state (C39)
$ws.Range("C39").Select()
